$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$ws = $wb.Worksheets.Item("Metadata")

# Insert a new row above "Description" (currently row 11) to make room for "Jurisdiction",
# pushing Description/Purpose/Copyright/Immutable down by one row.
$ws.Rows.Item(11).Insert()

# Match formatting of the surrounding property rows (style/border/wrap) for the new row.
$ws.Range("A12:B12").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)

# New Jurisdiction row (row 11), empty value
$ws.Cells.Item(11, 1).Value = "Jurisdiction"
$ws.Cells.Item(11, 2).Value = ""

# Bump the version + date metadata values
$ws.Cells.Item(3, 2).Value = "2.0.1-sd-202510-matchbox-patch"
$ws.Cells.Item(8, 2).Value = "2025-10-29T22:15:57+01:00"

# --- Rename the Include sheets ---
$wb.Worksheets.Item("Include from AddressUse").Name = "Include #0"
$wb.Worksheets.Item("Include from EntityNameUse").Name = "Include #1"
